# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Map of cell address -> new value, applied identically on both sheets.
$updates = @{
    "F2"  = 262
    "F3"  = 1343
    "F4"  = 150
    "F6"  = 227
    "F7"  = 97
    "F11" = 4529
    "F12" = 6797
    "F19" = 489
    "F21" = 57
    "F22" = 2701
    "F28" = 396
    "F30" = 35
    "F31" = 1624
    "F32" = 1020
    "F34" = 132
    "F35" = 80
    "F36" = 543
    "F39" = 89
    "F40" = 64
    "F42" = 11
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
